$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1. Title: A-E-C-P -> A-E-K-P
Replace-Text "A-E-C-P: Aseptá, Eksplorá, Konektá i Praktiká" "A-E-K-P: Aseptá, Eksplorá, Konektá i Praktiká"

# 2. Remove comma before "i permití" in the method description paragraph
Replace-Text "enfatisá prinsipionan sentral, i permití oportunidatnan" "enfatisá prinsipionan sentral i permití oportunidatnan"

# 3. "A-E-C-P ta para pa e 4 pasonan klave" -> "A-E-K-P ..."
Replace-Text "A-E-C-P ta para pa e 4 pasonan klave" "A-E-K-P ta para pa e 4 pasonan klave"

# 4. "Por apliká A-C-E-P na kualke momento" -> "Por apliká A-K-E-P na kualke momentu"
Replace-Text "Por apliká A-C-E-P na kualke momento durante" "Por apliká A-K-E-P na kualke momentu durante"

# 5. "duna ehèmpel di e A-E-C-P pa asina" -> "duna ehèmpel di e A-E-K-P pa asina"
Replace-Text "duna ehèmpel di e A-E-C-P pa asina" "duna ehèmpel di e A-E-K-P pa asina"

# 6. "Thank you for sharing." -> "Danki pa kompartí."
Replace-Text "Thank you for sharing. " "Danki pa kompartí. "

# 7. "Kon siguimentu di e ehèmpel di bo yu a laga bo sinti?" -> "Kon siguimentu di bo yu su guia a laga bo sinti?"
Replace-Text "Kon siguimentu di e ehèmpel di bo yu a laga bo sinti?" "Kon siguimentu di bo yu su guia a laga bo sinti?"

# 8. "Enkurashá mayornan pa traha nan mes konekshonnan" -> "... pa krea nan mes konekshonnan"
Replace-Text "Enkurashá mayornan pa traha nan mes konekshonnan" "Enkurashá mayornan pa krea nan mes konekshonnan"

# 9. "ku nos atenshon kompleto." -> "ku nos atenshon kompletu."
Replace-Text "ku nos atenshon kompleto. Esaki" "ku nos atenshon kompletu. Esaki"
